# "Generate Report for Handoff" -- refresh the localization-status report
# with the results of a new handoff run (new commit UUID / diff-hash and
# updated handoff timestamps), while leaving every hyperlink's underlying
# target URL exactly as it was (only the on-sheet display text and the
# cell content move to the new names/timestamps).

$wb = $excel.ActiveWorkbook

$oldName = "c55be2ed-1a84-48bd-ab73-9cd209c0f908"
$newName = "9d85a279-3a0d-4200-87b6-5e8786a6ac46"

$oldHash = "6947e14122ac0585d757c0a76c76ab30740313f2"
$newHash = "767b1805976f5542e88462388962767280d073d8"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/818aec78dcfbac96f9128897f01ce0909d818def/e2e/$oldName.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f863f5a758d550d70a4b3fbf1b7612b479af4edd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldName.$oldHash.zh-cn.xlf"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea414405f923d10d11d8e4316563bab8871504cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldName.$oldHash.de-de.xlf"

# ---- Sheet 1: Overview ----
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Hyperlinks.Delete()
$wsOverview.Range("A2").Value = "$newName.md"
$wsOverview.Range("D2").Value = "2016-56-18 05:56:25"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", "$newName.md") | Out-Null

# ---- Sheet 2: zh-cn ----
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("A2").Value = "$newName.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("D2").Value = "$newName.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-18 05:56:23"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", "$newName.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhUrl, "", "", "$newName.$newHash.zh-cn.xlf") | Out-Null

# ---- Sheet 3: de-de ----
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("A2").Value = "$newName.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("D2").Value = "$newName.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-18 05:56:25"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", "$newName.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deUrl, "", "", "$newName.$newHash.de-de.xlf") | Out-Null
